# Update the two-digit-division worksheet numbers in the single table.
# Each non-blank row (1, 5, 9, 13, 17) holds 5 division problems (columns 1-5).
# We address every changed cell directly by (row, column) and overwrite its
# run text, leaving the unchanged cell (row 17, column 3: "59÷5=") untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $table.Cell($row, $col).Range.Text = $text
}

# Row 1
Set-CellText $t 1 1 "33÷5="
Set-CellText $t 1 2 "57÷3="
Set-CellText $t 1 3 "66÷2="
Set-CellText $t 1 4 "11÷5="
Set-CellText $t 1 5 "23÷7="

# Row 5
Set-CellText $t 5 1 "72÷2="
Set-CellText $t 5 2 "18÷2="
Set-CellText $t 5 3 "35÷4="
Set-CellText $t 5 4 "31÷2="
Set-CellText $t 5 5 "29÷6="

# Row 9
Set-CellText $t 9 1 "25÷6="
Set-CellText $t 9 2 "96÷3="
Set-CellText $t 9 3 "97÷3="
Set-CellText $t 9 4 "38÷2="
Set-CellText $t 9 5 "86÷7="

# Row 13
Set-CellText $t 13 1 "85÷4="
Set-CellText $t 13 2 "27÷4="
Set-CellText $t 13 3 "83÷8="
Set-CellText $t 13 4 "74÷3="
Set-CellText $t 13 5 "12÷9="

# Row 17 (column 3, "59÷5=", is left unchanged per the diff)
Set-CellText $t 17 1 "24÷8="
Set-CellText $t 17 2 "39÷6="
Set-CellText $t 17 4 "52÷7="
Set-CellText $t 17 5 "65÷2="
